$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "women training pants",
    "women training suit full",
    "women training tights",
    "women travel shorts",
    "women underwear leggings",
    "women warmup pants",
    "women wetsuit pants",
    "women white capri",
    "women white compression shorts",
    "women wind pants",
    "women winter cycling pants",
    "women winter running gear",
    "women winter tights and leggings",
    "women winter tights thick",
    "women winter workout clothes",
    "women workout capris",
    "women workout pant",
    "women workout pants long",
    "women yoga capris",
    "women yoga fitness running gym stretch sports active leggings",
    "women’s capris",
    "women’s tights",
    "women’s wear",
    "womens 7 inch shorts",
    "womens active capri leggings",
    "womens active clothes",
    "womens active legging",
    "womens active leggings",
    "womens active pants",
    "womens active wear",
    "womens addias track pants",
    "womens adidas compression pants",
    "womens aparel",
    "womens apparal",
    "womens apparel clearance",
    "womens army leggings",
    "womens athleta",
    "womens athleta clothing",
    "womens athleta tops",
    "womens athletic capri",
    "womens athletic capri leggings",
    "womens athletic capri pants",
    "womens athletic capris",
    "womens athletic clothes",
    "womens athletic compression",
    "womens athletic gear",
    "womens athletic leggings black",
    "womens athletic leggings capri",
    "womens athletic leggings high waist",
    "womens athletic leggings yoga",
    "womens athletic pants",
    "womens athletic pants capri",
    "womens athletic pants tall",
    "womens athletic tights",
    "womens athletic wear",
    "womens ballet tights",
    "womens basketball",
    "womens basketball clothing",
    "womens basketball jersey",
    "womens basketball pants",
    "womens basketball short",
    "womens basketball shorts",
    "womens basketball shorts long",
    "womens best leggings",
    "womens best tights",
    "womens bicycle capris",
    "womens bike capris with padding",
    "womens bike tights",
    "womens bike wear",
    "womens bike with gears",
    "womens biking capris with padding",
    "womens biking gear",
    "womens biking shorts",
    "womens biking tights",
    "womens bill fold",
    "womens bills apparel",
    "womens black",
    "womens black and white stripes leggings",
    "womens black capri legging",
    "womens black capri leggings",
    "womens black capri pants",
    "womens black knee highs",
    "womens black nike shorts",
    "womens black pants",
    "womens black running tights",
    "womens black stocking",
    "womens black tights",
    "womens black workout leggings",
    "womens black yoga capris",
    "womens blue running shorts",
    "womens boxing apparel",
    "womens boxing gear",
    "womens boxing shorts",
    "womens bright blue tights",
    "womens brown leggins",
    "womens capri",
    "womens capri black",
    "womens capri black leggings",
    "womens capri leggings",
    "womens capri leggings plus",
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
